$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - F column updates ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 297
$ws1.Range("F3").Value = 62
$ws1.Range("F4").Value = 3634
$ws1.Range("F5").Value = 2227
$ws1.Range("F8").Value = 175
$ws1.Range("F9").Value = 92
$ws1.Range("F10").Value = 72
$ws1.Range("F11").Value = 1338
$ws1.Range("F13").Value = 1983
$ws1.Range("F14").Value = 144

# Sheet "全部类型" (All Types) - F column updates ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 297
$ws4.Range("F3").Value = 62
$ws4.Range("F4").Value = 3634
$ws4.Range("F5").Value = 2227
$ws4.Range("F9").Value = 175
$ws4.Range("F10").Value = 92
$ws4.Range("F11").Value = 72
$ws4.Range("F14").Value = 1338
$ws4.Range("F16").Value = 1983
$ws4.Range("F17").Value = 144
